$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 519.375
$ws.Range("I6").Value = 70
$ws.Range("J6").Value = 583.5714
$ws.Range("K6").Value = 210
$ws.Range("L6").Value = 1750.7142
$ws.Range("M6").Value = -98
$ws.Range("N6").Value = -1974.7142
$ws.Range("H15").Value = 719
$ws.Range("I15").Value = 719
$ws.Range("K15").Value = 2157
$ws.Range("M15").Value = -1988
$ws.Range("H18").Value = 1024
$ws.Range("I18").Value = 1024
$ws.Range("K18").Value = 1024
$ws.Range("M18").Value = -740
$ws.Range("H28").Value = 775
$ws.Range("I28").Value = 775
$ws.Range("K28").Value = 775
$ws.Range("M28").Value = -290
$ws.Range("H51").Value = 10400
$ws.Range("I51").Value = 12000
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = -11516
$ws.Range("N51").Value = -10968
$ws.Range("H62").Value = 5400
$ws.Range("I62").Value = 4250
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 4250
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -3626
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 5400
$ws.Range("I65").Value = 4250
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 21250
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -18130
$ws.Range("N65").Value = -56240
$ws.Range("H70").Value = 1642.1428
$ws.Range("I70").Value = 1498.3334
$ws.Range("J70").Value = 1750
$ws.Range("K70").Value = 4495.0002
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = -4225.0002
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 1642.1428
$ws.Range("I73").Value = 1498.3334
$ws.Range("J73").Value = 1750
$ws.Range("K73").Value = 4495.0002
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = -3559.0002
$ws.Range("N73").Value = -7122
$ws.Range("H80").Value = 1007.1429
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 810
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 2430
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -4426
$ws.Range("H83").Value = 1007.1429
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 810
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 7290
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -17274
$ws.Range("H88").Value = 3775.3333
$ws.Range("I88").Value = 4501
$ws.Range("J88").Value = 2324
$ws.Range("K88").Value = 4501
$ws.Range("L88").Value = 2324
$ws.Range("M88").Value = -4095
$ws.Range("N88").Value = -3136
$ws.Range("H91").Value = 3775.3333
$ws.Range("I91").Value = 4501
$ws.Range("J91").Value = 2324
$ws.Range("K91").Value = 4501
$ws.Range("L91").Value = 2324
$ws.Range("M91").Value = -3097
$ws.Range("N91").Value = -5132
$ws.Range("H138").Value = 1924.25
$ws.Range("I138").Value = 372.75
$ws.Range("J138").Value = 2700
$ws.Range("K138").Value = 1118.25
$ws.Range("L138").Value = 8100
$ws.Range("M138").Value = 4021.75
$ws.Range("N138").Value = -18380

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1918.25
$ws.Range("J88").Value = 1907.6666
$ws.Range("L88").Value = 1907.6666
$ws.Range("N88").Value = -2719.6666
$ws.Range("H91").Value = 1918.25
$ws.Range("J91").Value = 1907.6666
$ws.Range("L91").Value = 1907.6666
$ws.Range("N91").Value = -4715.6666
$ws.Range("H132").Value = 8075.5
$ws.Range("I132").Value = 2525
$ws.Range("J132").Value = 9925.666999999999
$ws.Range("K132").Value = 7575
$ws.Range("L132").Value = 29777.001
$ws.Range("M132").Value = -5045
$ws.Range("N132").Value = -34837.001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1668.75
$ws.Range("I86").Value = 1759.375
$ws.Range("K86").Value = 1759.375
$ws.Range("M86").Value = -636.375
$ws.Range("H89").Value = 1668.75
$ws.Range("I89").Value = 1759.375
$ws.Range("K89").Value = 8796.875
$ws.Range("M89").Value = -3180.875
$ws.Range("H105").Value = 1297.7142
$ws.Range("I105").Value = 1590
$ws.Range("K105").Value = 1590
$ws.Range("M105").Value = 157

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 30.4
$ws.Range("J7").Value = 52.5
$ws.Range("L7").Value = 52.5
$ws.Range("N7").Value = -278.5
$ws.Range("H36").Value = 7000
$ws.Range("I36").Value = 7000
$ws.Range("K36").Value = 7000
$ws.Range("M36").Value = -6612
$ws.Range("H40").Value = 7000
$ws.Range("I40").Value = 7000
$ws.Range("K40").Value = 7000
$ws.Range("M40").Value = -6840
$ws.Range("H132").Value = 3160
$ws.Range("I132").Value = 1825
$ws.Range("J132").Value = 4495
$ws.Range("K132").Value = 5475
$ws.Range("L132").Value = 13485
$ws.Range("M132").Value = -2945
$ws.Range("N132").Value = -18545
$ws.Range("H134").Value = 4200
$ws.Range("I134").Value = 5993
$ws.Range("J134").Value = 2407
$ws.Range("K134").Value = 17979
$ws.Range("L134").Value = 7221
$ws.Range("M134").Value = -15444
$ws.Range("N134").Value = -12291

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 74.833336
$ws.Range("I7").Value = 66
$ws.Range("J7").Value = 83.666664
$ws.Range("K7").Value = 198
$ws.Range("L7").Value = 250.999992
$ws.Range("M7").Value = -86
$ws.Range("N7").Value = -474.999992
$ws.Range("H11").Value = 500
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 1500
$ws.Range("N11").Value = -1780
$ws.Range("H23").Value = 339.77777
$ws.Range("I23").Value = 224.66667
$ws.Range("J23").Value = 397.33334
$ws.Range("K23").Value = 674.00001
$ws.Range("L23").Value = 1192.00002
$ws.Range("M23").Value = -439.00001
$ws.Range("N23").Value = -1662.00002
$ws.Range("H34").Value = 1460
$ws.Range("J34").Value = 1700
$ws.Range("L34").Value = 5100
$ws.Range("N34").Value = -5268
$ws.Range("H39").Value = 1395.1
$ws.Range("I39").Value = 391.66666
$ws.Range("J39").Value = 2900.25
$ws.Range("K39").Value = 1174.99998
$ws.Range("L39").Value = 8700.75
$ws.Range("M39").Value = -880.9999800000001
$ws.Range("N39").Value = -9288.75
$ws.Range("H55").Value = 863.3333
$ws.Range("I55").Value = 700
$ws.Range("K55").Value = 2100
$ws.Range("M55").Value = -1923
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 7500
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 12000
$ws.Range("K58").Value = 9000
$ws.Range("L58").Value = 36000
$ws.Range("M58").Value = -8872
$ws.Range("N58").Value = -36256
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6607.5454
$ws.Range("I70").Value = 6299.5
$ws.Range("J70").Value = 6783.5713
$ws.Range("K70").Value = 6299.5
$ws.Range("L70").Value = 6783.5713
$ws.Range("M70").Value = -6029.5
$ws.Range("N70").Value = -7323.5713
$ws.Range("H73").Value = 6607.5454
$ws.Range("I73").Value = 6299.5
$ws.Range("J73").Value = 6783.5713
$ws.Range("K73").Value = 6299.5
$ws.Range("L73").Value = 6783.5713
$ws.Range("M73").Value = -5363.5
$ws.Range("N73").Value = -8655.5713
$ws.Range("H107").Value = 226.85715
$ws.Range("I107").Value = 257.6
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 257.6
$ws.Range("L107").Value = 150
$ws.Range("M107").Value = 1662.4
$ws.Range("N107").Value = -3990
$ws.Range("H132").Value = 5550
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 993.125
$ws.Range("I46").Value = 990.8333
$ws.Range("K46").Value = 990.8333
$ws.Range("M46").Value = -802.8333
$ws.Range("H132").Value = 5168.4
$ws.Range("I132").Value = 5170.143
$ws.Range("K132").Value = 15510.429
$ws.Range("M132").Value = -12980.429

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 950
$ws.Range("I81").Value = 950
$ws.Range("K81").Value = 1900
$ws.Range("M81").Value = -839
$ws.Range("H84").Value = 950
$ws.Range("I84").Value = 950
$ws.Range("K84").Value = 9500
$ws.Range("M84").Value = -4196
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 1150
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3450
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -920
$ws.Range("N132").ClearContents()
